$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric need to be forced to Text so they
# keep the original "stringy" formatting (e.g. trailing zero in "1.00"),
# matching how this automation-generated sheet stores prices as text.
$textCells = @("D5", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D16", "D17", "D21", "D22", "D23", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D33", "D34", "D35", "D37", "D39", "D41", "D42", "D43", "D44", "D46", "D47", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '42.898.95'
$ws.Range("E2").Value = '  +3.64%  '
$ws.Range("D3").Value = '2.259.46'
$ws.Range("E3").Value = '  +3.20%  '
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").Value = '253.66'
$ws.Range("E5").Value = '  -0.49%  '
$ws.Range("E6").Value = '  +1.19%  '
$ws.Range("D7").Value = '72.14'
$ws.Range("E7").Value = '  +5.31%  '
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  -0.10%  '
$ws.Range("D9").Value = '0.647'
$ws.Range("D10").Value = '41.28'
$ws.Range("E10").Value = '  +9.09%  '
$ws.Range("D11").Value = '59.78'
$ws.Range("E11").Value = '  +1.29%  '
$ws.Range("D12").Value = '0.0967'
$ws.Range("E12").Value = '  +3.38%  '
$ws.Range("D13").Value = '7.41'
$ws.Range("E13").Value = '  +3.60%  '
$ws.Range("D14").Value = '0.104'
$ws.Range("E14").Value = '  +0.49%  '
$ws.Range("D15").Value = '2.596.39'
$ws.Range("E15").Value = '  +3.16%  '
$ws.Range("D16").Value = '0.886'
$ws.Range("E16").Value = '  +1.13%  '
$ws.Range("D17").Value = '14.82'
$ws.Range("E17").Value = '  +2.11%  '
$ws.Range("D18").Value = '2.252.98'
$ws.Range("E18").Value = '  +2.27%  '
$ws.Range("D19").Value = '42.846.31'
$ws.Range("E19").Value = '  +3.64%  '
$ws.Range("D20").Value = '0.0₃0977'
$ws.Range("E20").Value = '  +2.08%  '
$ws.Range("D21").Value = '6.25'
$ws.Range("E21").Value = '  +1.09%  '
$ws.Range("D22").Value = '73.17'
$ws.Range("E22").Value = '  +1.58%  '
$ws.Range("D23").Value = '236.29'
$ws.Range("E23").Value = '  +1.43%  '
$ws.Range("E24").Value = '  +4.46%  '
$ws.Range("D25").Value = '3.98'
$ws.Range("E25").Value = '  +0.91%  '
$ws.Range("D26").Value = '11.69'
$ws.Range("E26").Value = '  -0.54%  '
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  -0.10%  '
$ws.Range("D28").Value = '2.46'
$ws.Range("E28").Value = '  -2.62%  '
$ws.Range("D29").Value = '3.67'
$ws.Range("E29").Value = '  -1.46%  '
$ws.Range("D30").Value = '2.22'
$ws.Range("E30").Value = '  +2.39%  '
$ws.Range("D31").Value = '167.75'
$ws.Range("E31").Value = '  -0.69%  '
$ws.Range("E32").Value = '  +1.65%  '
$ws.Range("D33").Value = '0.129'
$ws.Range("E33").Value = '  +9.69%  '
$ws.Range("D34").Value = '6.17'
$ws.Range("E34").Value = '  +12.57%  '
$ws.Range("D35").Value = '0.0786'
$ws.Range("E35").Value = '  +3.78%  '
$ws.Range("E36").Value = '  +1.30%  '
$ws.Range("D37").Value = '28.93'
$ws.Range("E37").Value = '  +8.34%  '
$ws.Range("E38").Value = '  +1.62%  '
$ws.Range("D39").Value = '4.12'
$ws.Range("E39").Value = '  -0.83%  '
$ws.Range("E40").Value = '  +7.08%  '
$ws.Range("B41").Value = 'LidoDAOToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D41").Value = '2.30'
$ws.Range("E41").Value = '  +4.34%  '
$ws.Range("B42").Value = 'THORChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D42").Value = '5.98'
$ws.Range("E42").Value = '  +5.19%  '
$ws.Range("D43").Value = '12.37'
$ws.Range("E43").Value = '  -1.48%  '
$ws.Range("D44").Value = '64.21'
$ws.Range("E44").Value = '  -0.34%  '
$ws.Range("E45").Value = '  -2.10%  '
$ws.Range("B46").Value = 'FraxShare'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D46").Value = '9.02'
$ws.Range("E46").Value = '  +4.22%  '
$ws.Range("B47").Value = 'Algorand'
$ws.Range("C47").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D47").Value = '0.201'
$ws.Range("E47").Value = '  -0.22%  '
$ws.Range("E48").Value = '  +0.78%  '
$ws.Range("E49").Value = '  +3.18%  '
$ws.Range("E50").Value = '  -0.10%  '
$ws.Range("D51").Value = '4.40'
$ws.Range("E51").Value = '  +3.15%  '

foreach ($addr in $textCells) {
    $ws.Range($addr).ClearFormats()
}

Write-Output "applied 96 cell changes"
